# Rename header columns to use format-version suffixes instead of old/new,
# freeze the header row, and turn the data range into an Excel Table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row cells ------------------------------------------------
$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Freeze the header row ---------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Convert the used range into an Excel Table -------------------------------
# Stash the header row's current formatting (bold/fill/border/wrap) in a scratch
# cell so it can be restored after the table is created. Creating a ListObject on
# a range whose header already carries explicit formatting makes the host bake
# that formatting into a header dxf; stashing + clearing + restoring afterwards
# keeps the original styles.xml untouched (no new dxf, no new cell formats).
$headerRange = $ws.Range("A1:U1")
$ws.Range("A1").Copy() | Out-Null
$ws.Range("W1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$headerRange.ClearFormats()

$dataRange = $ws.Range("A1:U65")
$table = $ws.ListObjects.Add(1, $dataRange, 0, 1)
$table.Name = "Table1"
$table.TableStyle = ""

$ws.Range("W1").Copy() | Out-Null
$headerRange.PasteSpecial(-4122) | Out-Null      # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("W1").Clear() | Out-Null
